$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("API")

# Update the Base URL column (C) for all data rows: new Hungama POC host
$ws.Range("C2:C8").Value = "http://172.16.2.96:9090/task_manager/v1/"

# Split column D away from the D:E bestFit group and widen it (~13.43 chars,
# closest value the engine's pixel-quantized ColumnWidth setter can reach)
$ws.Columns.Item(4).ColumnWidth = 12.6

# Move the selection / active cell to C8 (also refreshes the frozen pane anchor)
$ws.Range("C8").Select()
$excel.ActiveWindow.FreezePanes = $true
